# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only column G ("K") values change for rows 2-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 2
    4  = 0
    5  = 5
    6  = 1
    7  = 0
    8  = 3
    9  = 2
    10 = 1
    11 = 2
    12 = 2
    13 = 1
    14 = 0
    15 = 6
    16 = 0
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 0
    23 = 2
    24 = 1
    25 = 2
    26 = 2
    27 = 2
    28 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
